$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 90 and 91 had their match data (columns F..V) swapped with each other.
# Columns A..E (index, pais, torneio, temporada, data_partida) are identical
# between the two rows already, so only F..V need to be rewritten.
# ---------------------------------------------------------------------------

# New contents for row 90 (previously row 91's match data)
$ws.Range("F90").Value = "IMT Novi Beograd"
$ws.Range("H90").Value = "Crvena zvezda"
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 8.15
$ws.Range("L90").Value = 24.2
$ws.Range("M90").Value = "28/10/2023 18:29"
$ws.Range("N90").Value = 5.6
$ws.Range("P90").Value = 9.529999999999999
$ws.Range("Q90").Value = "28/10/2023 18:29"
$ws.Range("R90").Value = 1.23
$ws.Range("T90").Value = 1.09
$ws.Range("U90").Value = "28/10/2023 18:21"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-crvena-zvezda/SjAgknkD/"

# New contents for row 91 (previously row 90's match data)
$ws.Range("F91").Value = "Radnicki Nis"
$ws.Range("H91").Value = "Sp. Subotica"
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = 1.7
$ws.Range("L91").Value = 1.66
$ws.Range("M91").Value = "28/10/2023 18:23"
$ws.Range("N91").Value = 3.42
$ws.Range("P91").Value = 3.65
$ws.Range("Q91").Value = "28/10/2023 18:23"
$ws.Range("R91").Value = 4.23
$ws.Range("T91").Value = 5.08
$ws.Range("U91").Value = "28/10/2023 18:23"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-spartak-subotica/2qDshl5f/"

# ---------------------------------------------------------------------------
# Append a brand-new row 103 with a new match (Radnicki Nis vs Novi Pazar).
# Formats are cloned from the last existing data row (102) so the index
# column keeps its bold/bordered style and the date column keeps its
# date-time number format, then values are written on top.
# ---------------------------------------------------------------------------

$ws.Range("A102").Copy()
$ws.Range("A103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E102").Copy()
$ws.Range("E103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "serbia"
$ws.Range("C103").Value = "super-liga"
$ws.Range("D103").Value = "2023-2024"
$ws.Range("E103").Value = 45240.70833333334
$ws.Range("F103").Value = "Radnicki Nis"
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = "Novi Pazar"
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 1.83
$ws.Range("K103").Value = "08/11/2023 17:13"
$ws.Range("L103").Value = 1.83
$ws.Range("M103").Value = "10/11/2023 16:57"
$ws.Range("N103").Value = 3.28
$ws.Range("O103").Value = "08/11/2023 17:13"
$ws.Range("P103").Value = 3.42
$ws.Range("Q103").Value = "10/11/2023 16:57"
$ws.Range("R103").Value = 3.77
$ws.Range("S103").Value = "08/11/2023 17:13"
$ws.Range("T103").Value = 4.31
$ws.Range("U103").Value = "10/11/2023 16:57"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-novi-pazar/UFGNPSQP/"
